# Update the Result/Date columns (A/B) on the various Katalon "Object
# repository" / test-result sheets to reflect a fresh test run, as recorded
# by the commit "Added Object repository for Pay As a guest page and ABp
# test suite".
#
# Each of these sheets has the layout:
#   A1=Result  B1=Date  C1=Notes  D1=Execute  ...
#   A2=<Pass/Fail>  B2=<run timestamp>  ...
# VerifyPasswordPolicy additionally stamps three rows (B2, B3, B4).

$wb = $excel.ActiveWorkbook

function Set-ResultRow($SheetName, $Result, $Timestamp) {
    $ws = $wb.Worksheets.Item($SheetName)
    $ws.Range("A2").Value = $Result
    $ws.Range("B2").Value = $Timestamp
}

Set-ResultRow "UIVerificationPendingBillsPage"  "Pass" "Thu Sep 25 14:01:13 IST 2025"
Set-ResultRow "CreateDeleteProfileOwner"        "Pass" "Thu Sep 25 14:10:39 IST 2025"
Set-ResultRow "CreateDeletePayer"               "Pass" "Thu Sep 25 14:09:47 IST 2025"
Set-ResultRow "VerifyUsernameLength"            "Pass" "Thu Sep 25 14:13:52 IST 2025"
Set-ResultRow "UiVerificationForAddUser"        "Pass" "Thu Sep 25 14:00:53 IST 2025"
Set-ResultRow "UiVerificationSPBillsLabel"      "Pass" "Wed Sep 24 17:22:49 IST 2025"
Set-ResultRow "UiVerificationSPIPDaily"         "Pass" "Thu Sep 25 14:04:21 IST 2025"
Set-ResultRow "UiVerificationSPInstallmentQuar" "Pass" "Thu Sep 25 14:09:03 IST 2025"
Set-ResultRow "UiVerificationSPInstallmentAnua" "Pass" "Thu Sep 25 14:08:18 IST 2025"
Set-ResultRow "UiVerificationSPIPDeferred"      "Pass" "Thu Sep 25 14:02:44 IST 2025"
Set-ResultRow "UiVerificationSPRecDeferred"     "Pass" "Thu Sep 25 14:04:58 IST 2025"
Set-ResultRow "UiVerificationSPAP"              "Fail" "Thu Sep 25 14:03:41 IST 2025"
Set-ResultRow "UiVerificationSPRecDaily"        "Fail" "Thu Sep 25 14:05:36 IST 2025"

# VerifyPasswordPolicy records three separate attempts (rows 2-4), each with
# its own timestamp in column B.
$wsPwd = $wb.Worksheets.Item("VerifyPasswordPolicy")
$wsPwd.Range("B2").Value = "Thu Sep 25 14:12:19 IST 2025"
$wsPwd.Range("B3").Value = "Thu Sep 25 14:12:53 IST 2025"
$wsPwd.Range("B4").Value = "Thu Sep 25 14:13:23 IST 2025"
